# "Fix CSV Example File"
#
# 1. The helper lookup sheet "Planilha4" (holds the lists used by the
#    data-validation dropdowns on "Relatório") is renamed to "Validação"
#    and hidden from the tab bar.
# 2. "Relatório Geral" ends up with two duplicate `_xlnm._FilterDatabase`
#    defined names (a leftover hidden/legacy one and the live/visible
#    one). Their ranges were swapped: the hidden one had been left
#    pointing at the live range ($A$2:$E$21) while the visible one
#    pointed at the stale range ($A$2:$E$22). Put each range back with
#    its correct visibility: hidden -> $A$2:$E$21, visible -> $A$2:$E$22.
# 3. The worksheet-level AutoFilter on "Relatório Geral" is shrunk from
#    A2:E22 to A2:E21 (row 22 has no data).
# 4. The data-validation list formulas on "Relatório" that pointed at
#    "Planilha4" are repointed at the renamed "Validação" sheet.
# 5. The lingering selection on "Relatório" (left over at A12) is reset
#    back to A2.

$wb = $excel.ActiveWorkbook

$wsValidacao   = $wb.Worksheets.Item(4)   # "Planilha4"
$wsRelGeral    = $wb.Worksheets.Item(2)   # "Relatório Geral"
$wsRelatorio   = $wb.Worksheets.Item(3)   # "Relatório"

# --- 1. rename + hide the validation helper sheet ------------------------
$wsValidacao.Name    = "Validação"
$wsValidacao.Visible = $xlSheetHidden

# --- 2. fix up the duplicated _xlnm._FilterDatabase defined names --------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Relatório Geral!_FilterDatabase") {
        if ($n.Visible) {
            $n.RefersTo = "='Relatório Geral'!`$A`$2:`$E`$22"
        } else {
            $n.RefersTo = "='Relatório Geral'!`$A`$2:`$E`$21"
        }
    }
}

# --- 3. resize the AutoFilter range on "Relatório Geral" -----------------
$wsRelGeral.AutoFilterMode = $false
$wsRelGeral.Range("A2:E21").AutoFilter()

# --- 4. repoint the data validation lists at the renamed sheet -----------
$wsRelatorio.Range("H2:H500").Validation.Formula1 = "Validação!`$A`$2:`$A`$3"
$wsRelatorio.Range("L2:L500").Validation.Formula1 = "Validação!`$B`$2:`$B`$3"
$wsRelatorio.Range("M2:M500").Validation.Formula1 = "Validação!`$C`$2:`$C`$3"
$wsRelatorio.Range("O2:O500").Validation.Formula1 = "Validação!`$D`$2:`$D`$6"
$wsRelatorio.Range("P2:P500").Validation.Formula1 = "Validação!`$E`$2:`$E`$4"

# --- 5. reset the stale selection on "Relatório" back to A2 --------------
$wsRelatorio.Activate()
$wsRelatorio.Range("A2").Select()
